$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 7: add description "GUI layout."
$ws.Range("D7").Value = "GUI layout."

# Row 8: update hours from 1.5 to 2
$ws.Range("C8").Value = 2

$wb.Application.Calculate()
